# Insert a new transaction row ("PAGO MOVISTAR") above the current row 1,
# pushing the existing five transactions down by one row, and correct the
# 'mo_monto' value on the "RETIRO CON LIBRETA" / 0031789350 row (the old
# row 5, now row 6) from 9.04 to 500.00.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at the top; everything currently on rows 1-6 shifts
#    down to rows 2-7 (formats/formulas move with the cells).
$ws.Rows.Item(1).Insert()

# 2. Populate the new row 1 with the "PAGO MOVISTAR" transaction.
$ws.Range("A1").Value = 41834
$ws.Range("B1").Value = "PAGO MOVISTAR"
$ws.Range("C1").Value = "D"
$ws.Range("D1").Value = "641925020506"
$ws.Range("E1").Value = "MATRIZ - QUITO"
$ws.Range("F1").Value = "3.00"
$ws.Range("G1").Value = "1.30"

$formula1 = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A1,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B1,"'', ''mo_tipo'' => ''",C1,"'', ''mo_documento'' => ''",D1,"'', ''mo_oficina'' => ''",E1,"'', ''mo_monto'' => ",F1,", ''mo_saldo'' => ",G1,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'
$ws.Range("H1").Formula = $formula1

# 3. Fix the 'mo_monto' amount on the "RETIRO CON LIBRETA" / 0031789350
#    transaction (now on row 6) — it should read 500.00, not 9.04.
$ws.Range("F6").Value = "500.00"

Write-Output "done"
